## Applies the commit:
##   1) Re-points the three "Table_0"-styled tables (slides 14-16) at the
##      built-in table style {0B17664E-9138-4132-A0C1-A07E427A95E5}.
##   2) Swaps the presentation's applied colour theme from the custom
##      "Integral / Red Violet" palette back to the stock "Office Theme"
##      palette (this is what the underlying theme1.xml/theme2.xml content
##      swap amounts to for the deck that's actually on screen - the
##      slide master's theme goes from Red Violet -> Office colours).

$p = $ppt.ActivePresentation

# --- 1) Table styles -------------------------------------------------
$oldTableStyle = "{D27B870B-1D52-44D6-A29C-AE1115C5CD75}"
$newTableStyle = "{0B17664E-9138-4132-A0C1-A07E427A95E5}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $sh = $s.Shapes.Item($i)
        if ($sh.HasTable) {
            $tbl = $sh.Table
            if ($tbl.Style -eq $oldTableStyle) {
                $tbl.ApplyStyle($newTableStyle)
            }
        }
    }
}

# --- 2) Theme colours: Red Violet -> Office --------------------------
# ppColorSchemeIndex order: 1 Background1(dk1) 2 Text1(lt1) 3 Background2(dk2)
# 4 Text2(lt2) 5 Accent1 .. 10 Accent6 11 Hyperlink 12 FollowedHyperlink
$officeColors = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

$s1 = $p.Slides.Item(1)
$tcs = $s1.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $tcs.Colors($i).RGB = $officeColors[$i - 1]
}
